$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 gets new claim data: new ID, new date, new time ---
# F2 (ID): old value "04104012736" is replaced by a brand-new value,
# which also makes the old shared string unused (it gets pruned on save).
$ws.Range("F6").Copy() | Out-Null
$ws.Range("F2").PasteSpecial(-4122) | Out-Null
$ws.Range("F2").Value = "04104015431"

# H2 (date) / I2 (time) get new values for row 2; use a leading apostrophe
# so the cells stay plain text with their original (quote-prefixed) styles.
$ws.Range("H2").Value = "'21/03/2022"
$ws.Range("I2").Value = "'12:00"

# --- Remove the hyperlink on C2 (row 2 no longer links out) ---
# C2 keeps its text value but loses the "Hipervinculo" style.
# (the engine's Hyperlinks.Delete()/Range.Hyperlinks wipe the whole sheet's
# hyperlinks, so rebuild C3:C9, then restore their original named style in
# one shot so the collection doesn't grow a stray duplicate per cell)
$url = "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/cc/ClaimCenter.do"
$ws.Hyperlinks.Delete()
foreach ($r in 3..9) {
    $ws.Hyperlinks.Add($ws.Range("C$r"), $url) | Out-Null
}
$ws.Range("C3:C9").Style = "Hipervínculo"
$ws.Range("C2").Style = "Normal"

# --- Move the active selection to F2 ---
$ws.Range("F2").Select() | Out-Null
